# Merge the three leading runs of the "Types can also be specified explicitly..."
# paragraph on the "Array Types" slide into a single run (text content is
# unchanged; only the run/XML structure collapses, matching the upstream
# OOXML diff which removed two superfluous <a:r> run splits).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$targetPrefix = "Types can also be specified explicitly. For Example: "

$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.StartsWith($targetPrefix)) {
        $sub = $tr.Characters($para.Start, $targetPrefix.Length)
        $sub.Text = $targetPrefix
        break
    }
}
